# BOM.xlsx update: found the D Pad part on Digikey instead of Sparkfun.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 is the "D Pad" component. Update vendor + link from Sparkfun to Digikey
# (same physical part, now sourced through Digikey).
$ws.Range("E11").Value = "Digikey"
$ws.Range("F11").Value = "https://www.digikey.com/en/products/detail/sparkfun-electronics/COM-26850/26266463"

# Leave the cursor where the author last edited (the new link cell).
$ws.Range("F13").Select() | Out-Null
